$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Fgf1"
$ws.Cells.Item(2, 3).Value = "Nrp1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 2.004760666666666
$ws.Cells.Item(2, 8).Value = 6.014282
$ws.Cells.Item(2, 9).Value = 0.1200698528618338
$ws.Cells.Item(2, 10).Value = 0.1200698528618338
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 102.8289443333334
$ws.Cells.Item(2, 14).Value = 308.486833
$ws.Cells.Item(2, 15).Value = 0.5559120396302444
$ws.Cells.Item(2, 16).Value = 0.5559120396302443
$ws.Cells.Item(2, 17).Value = 206.1474229943229
$ws.Cells.Item(2, 18).Value = 1855.326806948906
$ws.Cells.Item(2, 19).Value = 0.06674827680252535
$ws.Cells.Item(2, 20).Value = 0.06674827680252533
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Fgf1"
$ws.Cells.Item(3, 3).Value = "Nrp1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 2.004760666666666
$ws.Cells.Item(3, 8).Value = 6.014282
$ws.Cells.Item(3, 9).Value = 0.1200698528618338
$ws.Cells.Item(3, 10).Value = 0.1200698528618338
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 63.66262833333334
$ws.Cells.Item(3, 14).Value = 190.987885
$ws.Cells.Item(3, 15).Value = 0.3441717873742006
$ws.Cells.Item(3, 16).Value = 0.3441717873742006
$ws.Cells.Item(3, 17).Value = 127.6283332192856
$ws.Cells.Item(3, 18).Value = 1148.65499897357
$ws.Cells.Item(3, 19).Value = 0.0413246558692146
$ws.Cells.Item(3, 20).Value = 0.0413246558692146
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Fgf1"
$ws.Cells.Item(4, 3).Value = "Nrp1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 2.004760666666666
$ws.Cells.Item(4, 8).Value = 6.014282
$ws.Cells.Item(4, 9).Value = 0.1200698528618338
$ws.Cells.Item(4, 10).Value = 0.1200698528618338
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 18.481835
$ws.Cells.Item(4, 14).Value = 55.445505
$ws.Cells.Item(4, 15).Value = 0.09991617299555507
$ws.Cells.Item(4, 16).Value = 0.09991617299555505
$ws.Cells.Item(4, 17).Value = 37.05165585582333
$ws.Cells.Item(4, 18).Value = 333.46490270241
$ws.Cells.Item(4, 19).Value = 0.01199692019009383
$ws.Cells.Item(4, 20).Value = 0.01199692019009382
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Fgf1"
$ws.Cells.Item(5, 3).Value = "Nrp1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 9.409654999999999
$ws.Cells.Item(5, 8).Value = 28.228965
$ws.Cells.Item(5, 9).Value = 0.5635664696121425
$ws.Cells.Item(5, 10).Value = 0.5635664696121424
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 102.8289443333334
$ws.Cells.Item(5, 14).Value = 308.486833
$ws.Cells.Item(5, 15).Value = 0.5559120396302444
$ws.Cells.Item(5, 16).Value = 0.5559120396302443
$ws.Cells.Item(5, 17).Value = 967.5848901908718
$ws.Cells.Item(5, 18).Value = 8708.264011717845
$ws.Cells.Item(5, 19).Value = 0.3132933855893023
$ws.Cells.Item(5, 20).Value = 0.3132933855893021
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Fgf1"
$ws.Cells.Item(6, 3).Value = "Nrp1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 9.409654999999999
$ws.Cells.Item(6, 8).Value = 28.228965
$ws.Cells.Item(6, 9).Value = 0.5635664696121425
$ws.Cells.Item(6, 10).Value = 0.5635664696121424
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 63.66262833333334
$ws.Cells.Item(6, 14).Value = 190.987885
$ws.Cells.Item(6, 15).Value = 0.3441717873742006
$ws.Cells.Item(6, 16).Value = 0.3441717873742006
$ws.Cells.Item(6, 17).Value = 599.0433690098916
$ws.Cells.Item(6, 18).Value = 5391.390321089025
$ws.Cells.Item(6, 19).Value = 0.1939636791505792
$ws.Cells.Item(6, 20).Value = 0.1939636791505792
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Fgf1"
$ws.Cells.Item(7, 3).Value = "Nrp1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 9.409654999999999
$ws.Cells.Item(7, 8).Value = 28.228965
$ws.Cells.Item(7, 9).Value = 0.5635664696121425
$ws.Cells.Item(7, 10).Value = 0.5635664696121424
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 18.481835
$ws.Cells.Item(7, 14).Value = 55.445505
$ws.Cells.Item(7, 15).Value = 0.09991617299555507
$ws.Cells.Item(7, 16).Value = 0.09991617299555505
$ws.Cells.Item(7, 17).Value = 173.907691116925
$ws.Cells.Item(7, 18).Value = 1565.169220052325
$ws.Cells.Item(7, 19).Value = 0.05630940487226106
$ws.Cells.Item(7, 20).Value = 0.05630940487226104
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Fgf1"
$ws.Cells.Item(8, 3).Value = "Nrp1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 5.282203999999999
$ws.Cells.Item(8, 8).Value = 15.846612
$ws.Cells.Item(8, 9).Value = 0.3163636775260238
$ws.Cells.Item(8, 10).Value = 0.3163636775260238
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 102.8289443333334
$ws.Cells.Item(8, 14).Value = 308.486833
$ws.Cells.Item(8, 15).Value = 0.5559120396302444
$ws.Cells.Item(8, 16).Value = 0.5559120396302443
$ws.Cells.Item(8, 17).Value = 543.1634610733107
$ws.Cells.Item(8, 18).Value = 4888.471149659797
$ws.Cells.Item(8, 19).Value = 0.1758703772384168
$ws.Cells.Item(8, 20).Value = 0.1758703772384168
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Fgf1"
$ws.Cells.Item(9, 3).Value = "Nrp1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 5.282203999999999
$ws.Cells.Item(9, 8).Value = 15.846612
$ws.Cells.Item(9, 9).Value = 0.3163636775260238
$ws.Cells.Item(9, 10).Value = 0.3163636775260238
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 63.66262833333334
$ws.Cells.Item(9, 14).Value = 190.987885
$ws.Cells.Item(9, 15).Value = 0.3441717873742006
$ws.Cells.Item(9, 16).Value = 0.3441717873742006
$ws.Cells.Item(9, 17).Value = 336.2789900328467
$ws.Cells.Item(9, 18).Value = 3026.51091029562
$ws.Cells.Item(9, 19).Value = 0.1088834523544068
$ws.Cells.Item(9, 20).Value = 0.1088834523544068
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Fgf1"
$ws.Cells.Item(10, 3).Value = "Nrp1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 5.282203999999999
$ws.Cells.Item(10, 8).Value = 15.846612
$ws.Cells.Item(10, 9).Value = 0.3163636775260238
$ws.Cells.Item(10, 10).Value = 0.3163636775260238
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 18.481835
$ws.Cells.Item(10, 14).Value = 55.445505
$ws.Cells.Item(10, 15).Value = 0.09991617299555507
$ws.Cells.Item(10, 16).Value = 0.09991617299555505
$ws.Cells.Item(10, 17).Value = 97.62482276433998
$ws.Cells.Item(10, 18).Value = 878.6234048790599
$ws.Cells.Item(10, 19).Value = 0.03160984793320019
$ws.Cells.Item(10, 20).Value = 0.03160984793320019
